$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.379369
$ws.Range("H2").Value = 31.138107
$ws.Range("I2").Value = 0.01614698522449884
$ws.Range("J2").Value = 0.01614698522449883
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 297.8183156666666
$ws.Range("N2").Value = 893.454947
$ws.Range("O2").Value = 0.8852156413092672
$ws.Range("P2").Value = 0.8852156413092673
$ws.Range("Q2").Value = 3091.166193262814
$ws.Range("R2").Value = 27820.49573936533
$ws.Range("S2").Value = 0.014293563880716
$ws.Range("T2").Value = 0.014293563880716

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.379369
$ws.Range("H3").Value = 31.138107
$ws.Range("I3").Value = 0.01614698522449884
$ws.Range("J3").Value = 0.01614698522449883
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 24.34034433333333
$ws.Range("N3").Value = 73.021033
$ws.Range("O3").Value = 0.07234764413494278
$ws.Range("P3").Value = 0.0723476441349428
$ws.Range("Q3").Value = 252.6374154227257
$ws.Range("R3").Value = 2273.736738804531
$ws.Range("S3").Value = 0.001168196340874221
$ws.Range("T3").Value = 0.001168196340874221

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.379369
$ws.Range("H4").Value = 31.138107
$ws.Range("I4").Value = 0.01614698522449884
$ws.Range("J4").Value = 0.01614698522449883
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.277234
$ws.Range("N4").Value = 42.831702
$ws.Range("O4").Value = 0.04243671455578994
$ws.Range("P4").Value = 0.04243671455578994
$ws.Range("Q4").Value = 148.188679985346
$ws.Range("R4").Value = 1333.698119868114
$ws.Range("S4").Value = 0.0006852250029086148
$ws.Range("T4").Value = 0.0006852250029086147

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 604.0312093333333
$ws.Range("H5").Value = 1812.093628
$ws.Range("I5").Value = 0.9396797639857967
$ws.Range("J5").Value = 0.9396797639857967
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 297.8183156666666
$ws.Range("N5").Value = 893.454947
$ws.Range("O5").Value = 0.8852156413092672
$ws.Range("P5").Value = 0.8852156413092673
$ws.Range("Q5").Value = 179891.5573737531
$ws.Range("R5").Value = 1619024.016363778
$ws.Range("S5").Value = 0.8318192249020279
$ws.Range("T5").Value = 0.831819224902028

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 604.0312093333333
$ws.Range("H6").Value = 1812.093628
$ws.Range("I6").Value = 0.9396797639857967
$ws.Range("J6").Value = 0.9396797639857967
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.34034433333333
$ws.Range("N6").Value = 73.021033
$ws.Range("O6").Value = 0.07234764413494278
$ws.Range("P6").Value = 0.0723476441349428
$ws.Range("Q6").Value = 14702.32762325308
$ws.Range("R6").Value = 132320.9486092777
$ws.Range("S6").Value = 0.06798361716565143
$ws.Range("T6").Value = 0.06798361716565145

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 604.0312093333333
$ws.Range("H7").Value = 1812.093628
$ws.Range("I7").Value = 0.9396797639857967
$ws.Range("J7").Value = 0.9396797639857967
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.277234
$ws.Range("N7").Value = 42.831702
$ws.Range("O7").Value = 0.04243671455578994
$ws.Range("P7").Value = 0.04243671455578994
$ws.Range("Q7").Value = 8623.894918954984
$ws.Range("R7").Value = 77615.05427059486
$ws.Range("S7").Value = 0.03987692191811731
$ws.Range("T7").Value = 0.03987692191811731

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.39480333333333
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04417325078970442
$ws.Range("J8").Value = 0.04417325078970442
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 297.8183156666666
$ws.Range("N8").Value = 893.454947
$ws.Range("O8").Value = 0.8852156413092672
$ws.Range("P8").Value = 0.8852156413092673
$ws.Range("Q8").Value = 8456.492502419585
$ws.Range("R8").Value = 76108.43252177627
$ws.Range("S8").Value = 0.03910285252652329
$ws.Range("T8").Value = 0.0391028525265233

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.39480333333333
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04417325078970442
$ws.Range("J9").Value = 0.04417325078970442
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 24.34034433333333
$ws.Range("N9").Value = 73.021033
$ws.Range("O9").Value = 0.07234764413494278
$ws.Range("P9").Value = 0.0723476441349428
$ws.Range("Q9").Value = 691.1392904106144
$ws.Range("R9").Value = 6220.25361369553
$ws.Range("S9").Value = 0.003195830628417115
$ws.Range("T9").Value = 0.003195830628417116

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.39480333333333
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04417325078970442
$ws.Range("J10").Value = 0.04417325078970442
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.277234
$ws.Range("N10").Value = 42.831702
$ws.Range("O10").Value = 0.04243671455578994
$ws.Range("P10").Value = 0.04243671455578994
$ws.Range("Q10").Value = 405.39925157398
$ws.Range("R10").Value = 3648.59326416582
$ws.Range("S10").Value = 0.001874567634764009
$ws.Range("T10").Value = 0.001874567634764009

